$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.226.41"
$ws.Range("D3").Value = "1.571.34"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").Value = "'211.88"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D8").Value = "'22.09"
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D9").Value = "'0.249"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").Value = "'0.0870"
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").Value = "1.794.51"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").Value = "1.570.09"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "'0.521"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "27.261.43"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'62.37"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "'0.0" + [char]0x2083 + "0704"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "'216.52"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "'7.43"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").Value = "'154.02"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").Value = "'6.72"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("E30").Value = "  +2.41%  "
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").Value = "1.453.66"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("E34").Value = "  +2.18%  "
$ws.Range("E35").Value = "  +4.87%  "
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("D38").Value = "'0.0167"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "'0.535"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("E40").Value = "  +2.52%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "'64.72"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").Value = "1.707.98"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").Value = "'85.97"
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("D49").Value = "'0.0" + [char]0x2086 + "0105"
$ws.Range("E49").Value = "  +4.27%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "'0.0960"
$ws.Range("E51").Value = "  +0.55%  "
